# Revert "new changes in ops (ordercreation & orderpage & order form)"
#
# The prior commit had inserted two new columns ("Typist" / "Typist QC")
# right after "Assignee_QA", appended two more columns ("Status" / "Tier")
# at the end of the table, and appended a brand new data row (row 3) with
# a second record (Be18-002 / Blount / SIPL0102 / SIPL5317 / ...).
# This script undoes all of that, restoring the sheet to its original
# 13-column x 2-row ("A1:M2") shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the centered "AL / Autauga" formatting back to columns I:J -------
# Before the revert, that centered style lives on K2:L2 (after the two
# extra "Typist" columns shifted everything two slots to the right).
# Capture/copy it now, before we start rewriting cell values, and also
# grab the plain style (currently on E2) to restore on K2:L2 once the
# "State"/"County" values move back there.
$ws.Range("K2:L2").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E2").Copy()
$ws.Range("K2:L2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = 0

# --- Remove the extra data row (row 3: Be18-002 / Blount / ...) -----------
$ws.Rows.Item(3).Delete()

# --- Shift "Client".."Tier" two columns to the left, closing the gap left
#     by the removed "Typist" / "Typist QC" columns --------------------
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Lob"
$ws.Range("G1").Value = "Process"
$ws.Range("H1").Value = "Product Name"
$ws.Range("I1").Value = "State"
$ws.Range("J1").Value = "County"
$ws.Range("K1").Value = "Municipality"
$ws.Range("L1").Value = "Status"
$ws.Range("M1").Value = "Tier"

$ws.Range("E2").Value = "Baseline Title"
$ws.Range("F2").Value = "Title"
$ws.Range("G2").Value = "Search"
$ws.Range("H2").Value = "COS"
$ws.Range("I2").Value = "AL"
$ws.Range("J2").Value = "Autauga"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1)"

# --- Clear out the now-unused trailing columns (old Status/Tier slots) ----
$ws.Range("N1:O2").Clear()
